# LOQ4230.docx restructuring: the free-text "answer" paragraphs that sit
# under the Objetivos / Docente(s) Responsável(eis) / Programa resumido /
# Programa / Avaliação / Bibliografia headings got rotated into new homes
# (two independent cycles of content, discovered by diffing the before
# and after XML). Paragraph/run structure, styles and character
# formatting (bold labels, italics, line breaks) are unchanged -- only
# the wording inside each of 11 runs moves around. We therefore locate
# each source run once (fresh Find against the pristine wording) and
# overwrite it in place using absolute character ranges, walking from
# the end of the document towards the start so already-computed offsets
# for not-yet-processed runs stay valid while earlier ones shrink/grow.

$d = $word.ActiveDocument

function Get-Bounds($searchText) {
    $r = $d.Content.Duplicate
    $ok = $r.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found: $searchText"
    }
    return @{ Start = $r.Start; End = $r.End }
}

$fornecer = "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia de Produção nos diversos sistemas de produção da indústria. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."
$toProvide = "To provide an opportunity to apply the fundamental knowledge of Industrial Engineering in the various production systems of the industry. Complementation of general curricular training. Psychological and social adaptation of the student to his future professional activity"
$docente = "5840560 - Marco Antonio Carvalho Pereira"
$plano = "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."
$specificWorkPlan = "Specific Work Plan. Realization of the Internship. Final and / or partial report."
$participacao = "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."
$participation = "Participation of the student in the selective process of companies or in the academic sector. Internship carried out under the supervision of the School of Engineering of Lorena, through the Department of Chemical Engineering. The content will be established individually in the Work Plan between the Internship Supervisor and the tutor, as long as related to the areas of Industrial Engineering. Presentation of final report and / or partial reports about the activities carried out during the internship."
$supervisao = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$mf = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."
$naoSera = "Não será oferecida recuperação."
$aSerDefinida = "A ser definida com o orientador em função das atividades desenvolvidas no estágio."

# Capture the absolute Start/End of every one of the 11 runs while the
# document still has its original (pristine) wording, and remember the
# new text each one needs to end up with.
$slot1 = @{ Bounds = (Get-Bounds $fornecer);         New = $plano }
$slot2 = @{ Bounds = (Get-Bounds $toProvide);        New = $specificWorkPlan }
$slot3 = @{ Bounds = (Get-Bounds $docente);          New = $fornecer }
$slot4 = @{ Bounds = (Get-Bounds $plano);            New = $participacao }
$slot5 = @{ Bounds = (Get-Bounds $specificWorkPlan); New = $toProvide }
$slot6 = @{ Bounds = (Get-Bounds $participacao);     New = $supervisao }
$slot7 = @{ Bounds = (Get-Bounds $participation);    New = $participation }
$slot8 = @{ Bounds = (Get-Bounds $supervisao);       New = $mf }
$slot9 = @{ Bounds = (Get-Bounds $mf);               New = $naoSera }
$slot10 = @{ Bounds = (Get-Bounds $naoSera);         New = $aSerDefinida }
$slot11 = @{ Bounds = (Get-Bounds $aSerDefinida);    New = $docente }

# Document order (and thus ascending Start offset) is exactly
# slot1 .. slot11, so apply them back-to-front (slot11 first, slot1
# last) -- that way a write never shifts the offsets of a slot that
# still has to be processed.
$orderedHighToLow = @($slot11, $slot10, $slot9, $slot8, $slot7, $slot6, $slot5, $slot4, $slot3, $slot2, $slot1)

foreach ($slot in $orderedHighToLow) {
    $range = $d.Range($slot.Bounds.Start, $slot.Bounds.End)
    $range.Text = $slot.New
}

Write-Host "done"
